$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a tiny floating-point rounding difference in the existing row 27 timestamp
$ws.Range("A27").Value = 44340.78333358218

# Append the newly retrieved row of data
$ws.Range("A28").Value = 44341.78346584992
$ws.Range("B28").Value = 73996
$ws.Range("C28").Value = 62213
$ws.Range("D28").Value = 3260
$ws.Range("E28").Value = 2094
$ws.Range("F28").Value = 1477
$ws.Range("G28").Value = 19238
$ws.Range("H28").Value = 1339
$ws.Range("I28").Value = 833
$ws.Range("J28").Value = 196
